$d = $word.ActiveDocument

# 1) "no se tuvo" -> "no tuvo"
$d.Content.Find.Execute("no se tuvo mucho problema", $true, $false, $false, $false, $false, $true, 1, $false, "no tuvo mucho problema", 2) | Out-Null
